$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the second week's table (rows 12-16), columns A-C ---
# Column A: group member names
$ws.Range("A12").Value = "邱志鹏"
$ws.Range("A13").Value = "黄立根"
$ws.Range("A14").Value = "黄俊贤"
$ws.Range("A15").Value = "李达波"
$ws.Range("A16").Value = "冯德志"

# Column B: plan content
$ws.Range("B12").Value = "修改完善普通用户群组管理用例图，合并用例图                "
$ws.Range("B13").Value = "修改完善普通用户好友管理用例图"
$ws.Range("B14").Value = "修改完善普通用户个人信息管理用例图"
$ws.Range("B15").Value = "修改完善普通用户群管理员用例图"
$ws.Range("B16").Value = "修改完善管理员用例"

# Column C: completion status
$ws.Range("C12").Value = "完成"
$ws.Range("C13").Value = "完成"
$ws.Range("C14").Value = "完成"
$ws.Range("C15").Value = "完成"
$ws.Range("C16").Value = "完成"

# Columns B & C in rows 12-16 pick up the same font/border styling already
# used by the first table's body rows (row 3), while column A keeps the
# style it already had.
$ws.Range("B3:C3").Copy()
$ws.Range("B12:C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column width changes ---
$ws.Columns.Item(2).ColumnWidth = 89
$ws.Columns.Item(4).ColumnWidth = 40

# --- Selection change ---
$null = $ws.Range("C16").Select()
